# Merge a requerimientos no funcionales
#
# The "Matriz RACI" sheet's table (Tabla3) was re-sorted: instead of being
# sorted by the "Mares" column (C), it is now sorted alphabetically by the
# first column, "Actividad / Recurso" (A). The user then left the
# workbook with the "Matriz RACI" tab active/selected and a row from the
# freshly re-sorted table selected.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Matriz RACI")

# Re-sort Tabla3 by its first column ("Actividad / Recurso") ascending,
# replacing the previous sort-by-"Mares" (column C) order.
$lo = $ws3.ListObjects.Item("Tabla3")
$lo.Sort.SortFields.Clear() | Out-Null
$lo.Sort.SortFields.Add($ws3.Range("A2:A15")) | Out-Null
$lo.Sort.Header = 1
$lo.Sort.Apply() | Out-Null

# Leave "Matriz RACI" as the active sheet/tab, with a row of the
# freshly-sorted table selected.
$ws3.Activate() | Out-Null
$ws3.Range("A13:E13").Select() | Out-Null
